$d = $word.ActiveDocument

$d.Content.Find.Execute("613÷5=122, 3", $true, $false, $false, $false, $false, $true, 1, $false, "879÷6=146, 3", 2) | Out-Null
$d.Content.Find.Execute("186÷6=31, 0", $true, $false, $false, $false, $false, $true, 1, $false, "204÷2=102, 0", 2) | Out-Null
$d.Content.Find.Execute("322÷3=107, 1", $true, $false, $false, $false, $false, $true, 1, $false, "454÷5=90, 4", 2) | Out-Null
$d.Content.Find.Execute("115÷7=16, 3", $true, $false, $false, $false, $false, $true, 1, $false, "161÷9=17, 8", 2) | Out-Null
$d.Content.Find.Execute("897÷6=149, 3", $true, $false, $false, $false, $false, $true, 1, $false, "424÷5=84, 4", 2) | Out-Null
$d.Content.Find.Execute("922÷6=153, 4", $true, $false, $false, $false, $false, $true, 1, $false, "651÷7=93, 0", 2) | Out-Null
$d.Content.Find.Execute("298÷6=49, 4", $true, $false, $false, $false, $false, $true, 1, $false, "202÷8=25, 2", 2) | Out-Null
$d.Content.Find.Execute("189÷5=37, 4", $true, $false, $false, $false, $false, $true, 1, $false, "973÷9=108, 1", 2) | Out-Null
$d.Content.Find.Execute("798÷8=99, 6", $true, $false, $false, $false, $false, $true, 1, $false, "129÷2=64, 1", 2) | Out-Null
$d.Content.Find.Execute("858÷2=429, 0", $true, $false, $false, $false, $false, $true, 1, $false, "250÷9=27, 7", 2) | Out-Null
$d.Content.Find.Execute("272÷4=68, 0", $true, $false, $false, $false, $false, $true, 1, $false, "549÷7=78, 3", 2) | Out-Null
$d.Content.Find.Execute("596÷7=85, 1", $true, $false, $false, $false, $false, $true, 1, $false, "201÷5=40, 1", 2) | Out-Null
$d.Content.Find.Execute("696÷7=99, 3", $true, $false, $false, $false, $false, $true, 1, $false, "661÷5=132, 1", 2) | Out-Null
$d.Content.Find.Execute("648÷6=108, 0", $true, $false, $false, $false, $false, $true, 1, $false, "819÷4=204, 3", 2) | Out-Null
$d.Content.Find.Execute("403÷8=50, 3", $true, $false, $false, $false, $false, $true, 1, $false, "601÷8=75, 1", 2) | Out-Null
$d.Content.Find.Execute("181÷4=45, 1", $true, $false, $false, $false, $false, $true, 1, $false, "369÷8=46, 1", 2) | Out-Null
$d.Content.Find.Execute("604÷5=120, 4", $true, $false, $false, $false, $false, $true, 1, $false, "694÷6=115, 4", 2) | Out-Null
$d.Content.Find.Execute("311÷4=77, 3", $true, $false, $false, $false, $false, $true, 1, $false, "276÷3=92, 0", 2) | Out-Null
$d.Content.Find.Execute("747÷8=93, 3", $true, $false, $false, $false, $false, $true, 1, $false, "540÷7=77, 1", 2) | Out-Null
$d.Content.Find.Execute("338÷5=67, 3", $true, $false, $false, $false, $false, $true, 1, $false, "561÷7=80, 1", 2) | Out-Null
$d.Content.Find.Execute("874÷3=291, 1", $true, $false, $false, $false, $false, $true, 1, $false, "530÷6=88, 2", 2) | Out-Null
$d.Content.Find.Execute("566÷9=62, 8", $true, $false, $false, $false, $false, $true, 1, $false, "596÷9=66, 2", 2) | Out-Null
$d.Content.Find.Execute("929÷5=185, 4", $true, $false, $false, $false, $false, $true, 1, $false, "123÷6=20, 3", 2) | Out-Null
$d.Content.Find.Execute("719÷8=89, 7", $true, $false, $false, $false, $false, $true, 1, $false, "250÷8=31, 2", 2) | Out-Null
$d.Content.Find.Execute("528÷2=264, 0", $true, $false, $false, $false, $false, $true, 1, $false, "968÷9=107, 5", 2) | Out-Null
